$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: several columns (L and M) carry a Text ("@") number format. Typing
# a plain numeric literal into such a cell makes genuine Excel store it as
# text, which is NOT what the source data has (the existing cells in those
# columns hold real numbers even though the column is Text-formatted, since
# the sheet was originally populated programmatically). To reproduce a true
# numeric literal while keeping the cell's Text format, flip the number
# format to General just for the write, then restore it.
# ---------------------------------------------------------------------------
function Set-NumericValue {
    param($range, $value)
    $fmt = $range.NumberFormat
    $range.NumberFormat = "General"
    $range.Value = $value
    $range.NumberFormat = $fmt
}

# ---------------------------------------------------------------------------
# Row 413 - one more hospital death recorded (L413: 0 -> 1)
# ---------------------------------------------------------------------------
Set-NumericValue $ws.Range("L413") 1

# ---------------------------------------------------------------------------
# Row 419 - updated new-case count and a new hospital death
# ---------------------------------------------------------------------------
$ws.Range("C419").Value = 140
Set-NumericValue $ws.Range("L419") 1

# ---------------------------------------------------------------------------
# Row 420 - updated new-case count and a new hospital death
# ---------------------------------------------------------------------------
$ws.Range("C420").Value = 112
Set-NumericValue $ws.Range("L420") 1

# ---------------------------------------------------------------------------
# Row 421 - newly-entered day of data (was blank)
# ---------------------------------------------------------------------------
$ws.Range("C421").Value = 87
$ws.Range("E421").Value = 8
$ws.Range("F421").Value = 7
$ws.Range("G421").Value = 40
Set-NumericValue $ws.Range("L421") 0
Set-NumericValue $ws.Range("M421") 0

# ---------------------------------------------------------------------------
# Row 422 - newly-entered day of data (was blank)
# ---------------------------------------------------------------------------
$ws.Range("C422").Value = 13
$ws.Range("E422").Value = 8
$ws.Range("F422").Value = 7
$ws.Range("G422").Value = 31
Set-NumericValue $ws.Range("L422") 0
Set-NumericValue $ws.Range("M422") 0

# ---------------------------------------------------------------------------
# View state: the frozen header (rows 1:2 / col A) stays put; only the
# active selection on the data pane moves to L3:M422.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("L3:M422").Select()
